# Apply the "Data Analysis result commit 1" edit to the Fairness_TradeOffs
# sheet: every raw Likert-style answer in columns A-I (rows 2-117) is
# re-centered by subtracting 3 (so a 1-5 scale becomes a -2..2 scale).
# Empty cells (I7:I117, which hold no values) are left untouched.
# Finally, move the active cell selection from J10 to I3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fairness_TradeOffs")

$range = $ws.Range("A2:I117")
$values = $range.Value2

$rowCount = $values.GetLength(0)
$colCount = $values.GetLength(1)

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cellValue = $values[$r, $c]
        if ($cellValue -ne $null) {
            $values[$r, $c] = $cellValue - 3
        }
    }
}

$range.Value2 = $values

# Update the sheet's active cell / selection to I3 (was J10).
$ws.Activate()
$ws.Range("I3").Select() | Out-Null
